$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column C to fit the new "Test Parameters" text that will be added
# (matches the author's later best-fit column width)
$ws.Columns("C").ColumnWidth = 37.5

# Fill in the (previously empty) "Test Parameters" column C, in the order the
# values were actually typed in (rows 2-8, then 11-15, then 9-10, then 16-18)
$ws.Range("C2").Value = "Correct username/password"
$ws.Range("C3").Value = "Missing or incorrect password"
$ws.Range("C4").Value = "Correct username/password"
$ws.Range("C5").Value = "Missing or incorrect password"
$ws.Range("C6").Value = "All fields entered with information"
$ws.Range("C7").Value = "All fields entered with information"
$ws.Range("C8").Value = "Missing information on fields"

$ws.Range("C11").Value = "Status filter set to one of two statuses"
$ws.Range("C12").Value = "Assigned filter set to one of two statuses"
$ws.Range("C13").Value = "Assigned filter set to one of three statuses"
$ws.Range("C14").Value = "Button clicked for My tickets"
$ws.Range("C15").Value = "Button clicked for All tickets"

$ws.Range("C9").Value = "Tech logged in"
$ws.Range("C10").Value = "User logged in"

$ws.Range("C16").Value = "Search input filled out"
$ws.Range("C17").Value = "Ticket is cliked on home page"
$ws.Range("C18").Value = "Comment information filled out"

# Insert a new row 19 for the "Invalid/Missing inputs on Comment creation" test case,
# pushing the former rows 19-21 (Closing/Re-opening/Assigning ticket) down to 20-22
$ws.Rows("19").Insert()

$ws.Range("A19").Value = "Jared "
$ws.Range("B19").Value = "Invalid/Missing inputs on Comment creation"
$ws.Range("C19").Value = "Comment information missing fields"
$ws.Range("D19").Value = "Comment error should be displayed"

# Fill in new column C values for the rows that were shifted down
$ws.Range("C20").Value = "Close ticket button clicked"
$ws.Range("C21").Value = "Re-open ticket button clicked"
$ws.Range("C22").Value = "Assign ticket button clicked"

$ws.Range("C24").Select()
